# Further cleaning to metadata
#
# 1. Correct the libraryProtocol value used throughout column K (was "E7760", should be "E7420").
# 2. Give the libraryProtocol column (K2:K57) a refreshed font.
# 3. Replace the literal FALSE booleans in the roboticLibraryPrep column (L2:L57)
#    with an explicit =FALSE() formula.
# 4. Update the sheet's scroll/selection state to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 57

# --- 1 & 2: fix the value and restyle the libraryProtocol column ---
$kRange = $ws.Range("K2:K57")
$kRange.Value = "E7420"
$kRange.Font.Name = "Arial"
$kRange.Font.Size = 11
$kRange.Font.Color = 0

# --- 3: roboticLibraryPrep becomes a live formula instead of a hard-coded boolean ---
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("L$row").Formula = "=FALSE()"
}

# --- 4: leave the selection/scroll on the column that was just edited ---
[void]$ws.Range("K2:K57").Select()
